$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "schubert-winterreise_193"
$ws.Range("B2").Value = "schubert-winterreise_3"
$ws.Range("C2").Value = 0.0875
$ws.Range("D2").Value = "[['C:min', 'F:7/C', 'A#'], ['C/G', 'G:7', 'C']]"
$ws.Range("E2").Value = "[['C#:min', 'F#:7/A#', 'B:maj'], ['C#:maj', 'G#:7', 'C#:maj']]"
$ws.Range("F2").Value = "[(38.2, 41.06), (251.5, 255.38)]"
$ws.Range("G2").Value = "[(31.34, 36.62), (0.48, 8.48)]"
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()

$ws.Range("A3").Value = "schubert-winterreise_97"
$ws.Range("B3").Value = "schubert-winterreise_160"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "[['A:maj', 'D:maj', 'E:9/G#', 'A:maj', 'B:7/F#', 'E:min', 'A:(3,5,b7,b9)', 'D:maj', 'G:maj', 'D:maj/F#', 'A:min/E', 'B:maj/D#', 'C:maj/E', 'D:7/F#', 'G:maj', 'E:min7/G', 'E:(3,5,b7,b9)/G#', 'A:maj']]"
$ws.Range("E3").Value = "[['A:maj', 'D:maj', 'E:9/G#', 'A:maj', 'B:7/F#', 'E:min', 'A:(3,5,b7,b9)', 'D:maj', 'G:maj', 'D:maj/F#', 'A:min/E', 'B:maj/D#', 'C:maj/E', 'D:7/F#', 'G:maj', 'E:min7/G', 'E:(3,5,b7,b9)/G#', 'A:maj']]"
$ws.Range("F3").Value = "[(60.76, 82.8)]"
$ws.Range("G3").Value = "[(70.84, 95.2)]"
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()

$ws.Range("A4").Value = "schubert-winterreise_61"
$ws.Range("B4").Value = "schubert-winterreise_118"
$ws.Range("C4").Value = 0.07728085867620751
$ws.Range("D4").Value = "[['G:min/A#', 'A:hdim7/D#', 'D:7', 'G:min']]"
$ws.Range("E4").Value = "[['A:min', 'B:hdim7/D', 'E:7', 'A:min']]"
$ws.Range("F4").Value = "[(97.62, 113.78)]"
$ws.Range("G4").Value = "[(11.06, 17.28)]"
$ws.Range("H4").Value = "spotify:track:68YORkKP9uvlOQFMZZZwH5"
$ws.Range("I4").Value = "spotify:track:2qCvEz2hEb92VFATqVvrht"

$ws.Range("A5").Value = "isophonics_111"
$ws.Range("B5").Value = "schubert-winterreise_197"
$ws.Range("C5").Value = 0.2015810276679842
$ws.Range("D5").Value = "[['C', 'G:7', 'C']]"
$ws.Range("E5").Value = "[['F#:maj/A#', 'C#:7', 'F#:maj']]"
$ws.Range("F5").Value = "[(70.519024, 75.720294)]"
$ws.Range("G5").Value = "[(17.48, 19.72)]"
$ws.Range("H5").ClearContents()
$ws.Range("I5").Value = "spotify:track:4lrfYSnZmpXdCWuWqVo8L0"

$ws.Range("A6").Value = "isophonics_165"
$ws.Range("B6").Value = "isophonics_212"
$ws.Range("C6").Value = 0.2509803921568627
$ws.Range("D6").Value = "[['G', 'A', 'D', 'A', 'D', 'A', 'D', 'A']]"
$ws.Range("E6").Value = "[['C', 'D', 'G', 'D', 'G', 'D', 'G', 'D']]"
$ws.Range("F6").Value = "[(40.281383, 49.731904)]"
$ws.Range("G6").Value = "[(46.491101, 61.200942)]"
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()

$ws.Range("A7").Value = "schubert-winterreise_153"
$ws.Range("B7").Value = "schubert-winterreise_57"
$ws.Range("C7").Value = 0.7142857142857143
$ws.Range("D7").Value = "[['D#:maj/A#', 'A#:(3,5,b7,b9)', 'D#:maj']]"
$ws.Range("E7").Value = "[['A#:maj', 'F:(3,5,b7,b9)', 'A#:maj']]"
$ws.Range("F7").Value = "[(54.78, 57.52)]"
$ws.Range("G7").Value = "[(23.16, 26.04)]"
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()

$ws.Range("A8").Value = "jaah_41"
$ws.Range("B8").Value = "jaah_85"
$ws.Range("C8").Value = 0.08221524407650585
$ws.Range("D8").Value = "[['F', 'F', 'F', 'F']]"
$ws.Range("E8").Value = "[['Ab', 'Ab', 'Ab', 'Ab']]"
$ws.Range("F8").Value = "[(133.33, 138.07)]"
$ws.Range("G8").Value = "[(83.46, 88.05)]"
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()

$ws.Range("A9").Value = "schubert-winterreise_154"
$ws.Range("B9").Value = "schubert-winterreise_157"
$ws.Range("C9").Value = 0.2687747035573123
$ws.Range("D9").Value = "[['E:7', 'A:maj', 'A:maj', 'E:7/G#']]"
$ws.Range("E9").Value = "[['C:7', 'F:maj', 'F:maj/A', 'C:7']]"
$ws.Range("F9").Value = "[(9.24, 17.72)]"
$ws.Range("G9").Value = "[(17.36, 20.12)]"
$ws.Range("H9").Value = "spotify:track:0XfunCHFEeQnzm4NaY8rJr"
$ws.Range("I9").Value = "spotify:track:4lrfYSnZmpXdCWuWqVo8L0"

$ws.Range("A10").Value = "schubert-winterreise_34"
$ws.Range("B10").Value = "schubert-winterreise_195"
$ws.Range("C10").Value = 0.09642857142857142
$ws.Range("D10").Value = "[['C:min/G', 'G', 'C:min']]"
$ws.Range("E10").Value = "[['F:min/C', 'C', 'F:min/C']]"
$ws.Range("F10").Value = "[(11.82, 16.4)]"
$ws.Range("G10").Value = "[(44.16, 48.16)]"
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()

$ws.Range("A11").Value = "schubert-winterreise_193"
$ws.Range("B11").Value = "jaah_79"
$ws.Range("C11").Value = 0.03709677419354838
$ws.Range("D11").Value = "[['C', 'C/G', 'G:7']]"
$ws.Range("E11").Value = "[['Bb', 'Bb', 'F:7']]"
$ws.Range("F11").Value = "[(250.3, 252.76)]"
$ws.Range("G11").Value = "[(7.59, 10.86)]"
$ws.Range("H11").ClearContents()
$ws.Range("I11").ClearContents()

$ws.Range("A12").Value = "schubert-winterreise_26"
$ws.Range("B12").Value = "schubert-winterreise_169"
$ws.Range("C12").Value = 0.2666666666666667
$ws.Range("D12").Value = "[['A#:min', 'F:min', 'C:7', 'F:min']]"
$ws.Range("E12").Value = "[['C:min', 'G:min', 'D:7', 'G:min']]"
$ws.Range("F12").Value = "[(8.18, 14.18)]"
$ws.Range("G12").Value = "[(17.86, 28.3)]"
$ws.Range("H12").Value = "spotify:track:1nvxQGWCnikMK7a4HYQvSx"
$ws.Range("I12").Value = "spotify:track:3OD2uwEUQKg0WyW9Lewata"

$ws.Range("A13").Value = "schubert-winterreise_151"
$ws.Range("B13").Value = "isophonics_216"
$ws.Range("C13").Value = 0.3483870967741935
$ws.Range("D13").Value = "[['C:maj/G', 'F:maj', 'C:maj/G']]"
$ws.Range("E13").Value = "[['A', 'D', 'A']]"
$ws.Range("F13").Value = "[(117.54, 121.8)]"
$ws.Range("G13").Value = "[(28.009637, 33.826235)]"
$ws.Range("H13").ClearContents()
$ws.Range("I13").ClearContents()

$ws.Range("A14").Value = "schubert-winterreise_119"
$ws.Range("B14").Value = "schubert-winterreise_51"
$ws.Range("C14").Value = 0.1916666666666667
$ws.Range("D14").Value = "[['G:maj', 'D:7', 'G:maj']]"
$ws.Range("E14").Value = "[['C/G', 'G:7', 'C']]"
$ws.Range("F14").Value = "[(7.6, 17.16)]"
$ws.Range("G14").Value = "[(254.26, 257.98)]"
$ws.Range("H14").ClearContents()
$ws.Range("I14").ClearContents()

$ws.Range("A15").Value = "isophonics_277"
$ws.Range("B15").Value = "isophonics_57"
$ws.Range("C15").Value = 0.08947368421052632
$ws.Range("D15").Value = "[['A', 'E', 'A'], ['A', 'E', 'D'], ['E', 'A', 'E']]"
$ws.Range("E15").Value = "[['C', 'G', 'C'], ['C', 'G', 'F'], ['G', 'C', 'G']]"
$ws.Range("F15").Value = "[(31.840929, 37.216349), (27.452358, 31.840929), (32.920657, 39.410634)]"
$ws.Range("G15").Value = "[(136.604739, 139.228594), (5.017256, 9.231677), (4.155032, 7.525011)]"
$ws.Range("H15").Value = "spotify:track:2RnPATK99oGOZygnD2GTO6"
$ws.Range("I15").ClearContents()

$ws.Range("A16").Value = "isophonics_295"
$ws.Range("B16").Value = "schubert-winterreise_14"
$ws.Range("C16").Value = 0.1666666666666667
$ws.Range("D16").Value = "[['C/3', 'F', 'C/3']]"
$ws.Range("E16").Value = "[['D:maj', 'G:maj', 'D:maj/F#']]"
$ws.Range("F16").Value = "[(76.585782, 78.06481)]"
$ws.Range("G16").Value = "[(83.34, 87.64)]"
$ws.Range("H16").ClearContents()
$ws.Range("I16").ClearContents()

$ws.Range("A17").Value = "schubert-winterreise_37"
$ws.Range("B17").Value = "jaah_87"
$ws.Range("C17").Value = 0.1430860805860806
$ws.Range("D17").Value = "[['F:min/C', 'C', 'F:min/C', 'C', 'F:min']]"
$ws.Range("E17").Value = "[['F:min', 'C', 'F:min', 'C', 'F:min']]"
$ws.Range("F17").Value = "[(45.58, 52.16)]"
$ws.Range("G17").Value = "[(0.6, 8.74)]"
$ws.Range("H17").ClearContents()
$ws.Range("I17").ClearContents()
